$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card11")

# Fill previously-empty inline-string cells (rows 2-14) with the literal text "nan",
# matching the re-export pattern used when the source sheet was regenerated.
$nanRefs = "D2","E2","F2","G2","H2","I2","J2","K2","L2","M2","N2","O2","P2","D3","G3","H3","I3","J3","K3","M3","N3","O3","P3","D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","N4","O4","P4","D5","H5","I5","J5","K5","M5","N5","O5","P5","E6","F6","G6","I6","J6","K6","M6","N6","O6","P6","E7","G7","H7","J7","M7","N7","O7","P7","E8","H8","I8","J8","K8","M8","N8","O8","P8","F9","G9","H9","I9","J9","K9","M9","N9","P9","E10","F10","G10","H10","I10","J10","K10","P10","D11","E11","F11","G11","H11","I11","J11","K11","L11","M11","N11","O11","P11","D12","E12","F12","G12","H12","I12","J12","K12","L12","M12","N12","O12","P12","D13","E13","F13","G13","H13","I13","J13","K13","L13","M13","N13","O13","P13","B14","C14","D14","E14","F14","G14","H14","I14","J14","K14","P14"

foreach ($ref in $nanRefs) {
    $ws.Range($ref).Value = "nan"
}

# Append the new service event as row 15. Force text format on the
# numeric-/date-looking values so they stay literal strings (matching the
# rest of the sheet, which stores "Date" and "card" columns as text).
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "11"
$ws.Range("L15").NumberFormat = "@"
$ws.Range("L15").Value = "9/9/2025"
$ws.Range("M15").Value = "قطع سير كويلر مسنن 1270"
$ws.Range("N15").Value = "تم تغير سير 1270(ميجا فلكس)"
$ws.Range("O15").Value = "فني"
